$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 4587  # was 4584
$ws.Range("F4").Value = 3704  # was 3703
$ws.Range("F6").Value = 173  # was 172
$ws.Range("F8").Value = 387  # was 384
$ws.Range("F9").Value = 387  # was 385
$ws.Range("F10").Value = 2618  # was 2614
$ws.Range("F18").Value = 73  # was 72
$ws.Range("F19").Value = 10829  # was 10810
$ws.Range("F20").Value = 6219  # was 6211
$ws.Range("F24").Value = 226  # was 222
$ws.Range("F29").Value = 207  # was 205
$ws.Range("F31").Value = 3583  # was 3581
$ws.Range("F33").Value = 974  # was 973
$ws.Range("F34").Value = 487  # was 486
$ws.Range("F36").Value = 287  # was 286
$ws.Range("F38").Value = 264  # was 263
$ws.Range("F39").Value = 4894  # was 4893
$ws.Range("F41").Value = 1172  # was 1169

$ws = $wb.Worksheets.Item(2)
$ws.Range("F13").Value = 3627  # was 3626

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8896  # was 8893
$ws.Range("F3").Value = 455  # was 454
$ws.Range("F4").Value = 1704  # was 1700

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 455  # was 454
$ws.Range("F3").Value = 1704  # was 1700
$ws.Range("F4").Value = 4587  # was 4584
$ws.Range("F5").Value = 3704  # was 3703
$ws.Range("F9").Value = 387  # was 385
$ws.Range("F10").Value = 2618  # was 2614
$ws.Range("F20").Value = 10829  # was 10810
$ws.Range("F21").Value = 3627  # was 3626
$ws.Range("F26").Value = 226  # was 222
$ws.Range("F30").Value = 207  # was 205
$ws.Range("F32").Value = 3583  # was 3581
$ws.Range("F34").Value = 974  # was 973
$ws.Range("F36").Value = 287  # was 286
$ws.Range("F40").Value = 264  # was 263
$ws.Range("F41").Value = 4894  # was 4893
$ws.Range("F43").Value = 1172  # was 1169
